# Adds the 2018-04-11 work-record entry after the existing 2018-04-10 entry.
#
# The new entry's paragraph text/runs (with their exact w:rFonts eastAsia hints,
# split runs, the literal tab character, and the trailing bookmark) are easiest
# to reproduce faithfully via Range.InsertXML with a WordprocessingML fragment,
# rather than rebuilding every run/formatting toggle through the object model.
$d = $word.ActiveDocument

# Sanity-check: the paragraph we are about to replace/extend should be the
# "2018-04-10" entry's body paragraph (the last paragraph in the document).
$targetPara = $d.Paragraphs.Last
$checkRange = $d.Content
$checkRange.Find.Execute("2018年4月10日 周二") | Out-Null
if ($checkRange.Start -ge $targetPara.Range.Start) {
    Write-Host "WARNING: anchor paragraph not found before the last paragraph"
}

$r = $targetPara.Range
# Exclude the trailing paragraph mark from the replaced range so InsertXML
# only swaps in new paragraphs without leaving a stray empty one behind.
$r.End = $r.End - 1

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00DE3912" w:rsidRPr="00E60E36" w:rsidRDefault="00DE3912" w:rsidP="00F16E84"><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>昨晚配置好了环境后，今天过来主要进行代码调试。用U盘将之前已经写好的代码拷到服务器上，进行第一次make尝试。根据学长的说法，由于我的pmfs采用内核模块方式，因此每次只需要在我自己的mypmfs中进行编译和安装模块就可以了。上午第一次make自己的代码，根据调试信息改了不少的问题，主要有：1）pmfs_blockp_alloc()函数定义的参数为空，但是要写void在里面；2）balloc.c中的BUG_ON;3)一些指针初始化为NULL；4）在pmfs</w:t></w:r><w:r><w:t>.h</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>中添加自己定义的函数的声明；5）添加kmem</w:t></w:r><w:r><w:t>_cache_alloc</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>相关的函数，主要有create，destroy等，修改了init_pmfs_fs()函数。在改完调试信息里面能看到的一些bug后，我尝试insmod和mount命令，但是在mount的时候出了问题，进程直接被kill了。通过dmesg查看日志信息，发现报错是在我的自定义函数pmfs_blockp_alloc()中，后来经过仔细分析源码，发现是我自己定义的pmfs_blockp_cachep没有进行初始化，于是仿照pmfs_blocknode_cache相关的几个函数，进行初始化，解决了这个问题</w:t></w:r><w:r w:rsidR="00DE77FE"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>，于是也成功的将pmfs挂载上去了。但是后面在用postmark进行测试的时候，又出现了新的问题</w:t></w:r><w:r w:rsidR="00E60E36"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>，报错信息比较模糊，通过dmesg查看日志发现错误定位在__pmfs_free_block()函数中，并且有精确定位balloc.c第144行，查看是BUG_ON(list</w:t></w:r><w:r w:rsidR="00E60E36"><w:t>_empty(head))</w:t></w:r><w:r w:rsidR="00E60E36"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>，而这个head之前初始化为sbi-&gt;block_inuse_head，是空的，猜想可能是这里出发了内核的预警，中断了postmark的操作。并且导致无法umount。</w:t></w:r><w:r w:rsidR="00B74B25"><w:tab/></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2018年4月11日 周三</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>今天主要是解决postmark无法完成测试的情况，在使用配置文件的情况下执行./</w:t></w:r><w:r><w:t>postmark pm.cfg</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>时，postmark程序会直接被kill掉，用dmesg查看日志信息发现错误是由于引用空指针触发了内核的BUG，导致程序崩溃，并且postmark崩溃，导致pmfs文件系统无法从卸载点umount，pmfs模块也无法用rmmod指令移除，在网上查看了大量的资料后，发现pmfs文件系统可以用umount</w:t></w:r><w:r><w:t xml:space="preserve"> –fl /mnt/pmfs</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>指令强制卸载，但是pmfs模块的卸载却遇到了很大的问题，在移除pmfs文件系统后，用lsmod命令查看pmfs模块，发现引用计数为1，这也就是pmfs模块无法卸载的原因（提示module is in use），根据网上的资料，需要另外写一个force</w:t></w:r><w:r><w:t>_rmmod</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>模块，来将pmfs模块的引用计数清零，但是我在尝试这种办法的时候，发现force_rmmod被插上之后，该模块也会直接崩溃，并且无法使用reboot命令重启机器（这是最严重的！），</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>并且用不同的方法实现force</w:t></w:r><w:r><w:t>_rmmod</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>模块有不同的效果，在使用complete函数时，甚至会导致机器直接重启（算是一种间接重启的方法）。后来放弃了模块的卸载，每次都是重启后再编译测试（没办法的办法）。自己写了个small</w:t></w:r><w:r><w:t>.cfg</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>配置文件，测试10个事务，并且在自己写的几个函数中添加了大量的debug信息，最终才定位到是自己在写pmfs_new_block()时，有一种情况中，没有将二级索引节点Pmfs_blockp与一级索引pmfs_blocknode相互连接起来，导致发生了空指针引用，解决完这个问题之后，</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:lastRenderedPageBreak/><w:t>postmark程序才真正能够完全运行测试，但是新的问题是，测试结果竟然与优化前的方案一模一样！性能既没有降低，也没有提升，一番考虑后认为可能是由于postmark的序列化操作，以及配置文件中设置的文件size过小导致这种结果，于是重新设置了新的配置文件my.cfg，并重新分配了4G内存给pmfs，测试后仍然是一样的，这让我不得不怀疑postmark程序并不能测试出分配器的优化，需要寻求其他测试工具的辅助。</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r.InsertXML($xml)

Write-Host "Paragraphs now: $($d.Paragraphs.Count)"
